$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Real Word's InlineShape object has no writable "Name" of its own, so the
# only way to relabel the drawing's <wp:docPr name="..."> is to briefly
# promote the inline picture to a floating Shape (which does expose
# .Name), rename it there, then convert it back to an inline picture in
# place (same position/formatting, no visible change to the user).
function Rename-HeaderFooterPicture($story, $newName) {
    $ils = $story.Range.InlineShapes.Item(1)
    $shape = $ils.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

# Headers.Item(1) is the default (odd-page) header, Headers.Item(2) is the
# first-page header -- both carry the Pearson/BTEC logo picture currently
# named "image2.jpg"; rename both to "image1.jpg".
Rename-HeaderFooterPicture $sec.Headers.Item(1) "image1.jpg"
Rename-HeaderFooterPicture $sec.Headers.Item(2) "image1.jpg"

# Footers.Item(1) is the default (odd-page) footer, Footers.Item(2) is the
# first-page footer -- both carry the Pearson logo picture currently named
# "image1.png"; rename both to "image2.png".
Rename-HeaderFooterPicture $sec.Footers.Item(1) "image2.png"
Rename-HeaderFooterPicture $sec.Footers.Item(2) "image2.png"

Write-Output "Renamed header/footer logo pictures"
